$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to be treated as text so Excel doesn't auto-convert
    # decimal-looking strings (e.g. "594.70") into numbers, then restore
    # the cell's original (default) style so no stray formatting is left
    # behind.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "67.404.06"
$ws.Range("E2").Value = "  -1.39%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.755.28"
$ws.Range("E3").Value = "  -1.96%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
Set-TextValue "D5" "594.70"
$ws.Range("E5").Value = "  -1.08%  "

# Row 6 - Solana
Set-TextValue "D6" "170.05"
$ws.Range("E6").Value = "  +0.06%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.753.98"
$ws.Range("E7").Value = "  -1.98%  "

# Row 9 - XRP
Set-TextValue "D9" "0.524"
$ws.Range("E9").Value = "  -0.50%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.74%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +0.37%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.84%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +4.68%  "

# Row 14 - Avalanche
Set-TextValue "D14" "36.53"
$ws.Range("E14").Value = "  -1.53%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.388.88"
$ws.Range("E15").Value = "  -1.90%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.757.72"
$ws.Range("E16").Value = "  -1.89%  "

# Row 17 - Chainlink
Set-TextValue "D17" "18.58"
$ws.Range("E17").Value = "  +0.51%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "67.430.00"
$ws.Range("E18").Value = "  -1.29%  "

# Row 19 - Polkadot
Set-TextValue "D19" "7.19"
$ws.Range("E19").Value = "  -3.05%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  +1.01%  "

# Row 21 - Uniswap
Set-TextValue "D21" "10.52"
$ws.Range("E21").Value = "  -5.41%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "467.15"
$ws.Range("E22").Value = "  -0.67%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.718"
$ws.Range("E23").Value = "  -2.16%  "

# Row 24 - Litecoin
Set-TextValue "D24" "83.90"
$ws.Range("E24").Value = "  +0.94%  "

# Row 25 - PEPE
Set-TextValue "D25" "0.0000146"
$ws.Range("E25").Value = "  -8.88%  "

# Row 26 - Fetch.AI
Set-TextValue "D26" "2.21"
$ws.Range("E26").Value = "  -1.12%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "12.14"
$ws.Range("E27").Value = "  -0.25%  "

# Row 28 - RenderToken
Set-TextValue "D28" "10.32"
$ws.Range("E28").Value = "  +2.99%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  -0.16%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -2.00%  "

# Row 31 - WrappedeETH
Set-TextValue "D31" "3.909.25"
$ws.Range("E31").Value = "  -1.77%  "

# Row 32 - NEARProtocol
Set-TextValue "D32" "7.65"
$ws.Range("E32").Value = "  -0.64%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "30.55"
$ws.Range("E33").Value = "  -3.24%  "

# Row 34 - ImmutableX
$ws.Range("E34").Value = "  -3.57%  "

# Row 35 - Aptos
Set-TextValue "D35" "9.12"
$ws.Range("E35").Value = "  -3.41%  "

# Row 36 - RenzoRestakedETH
Set-TextValue "D36" "3.721.65"
$ws.Range("E36").Value = "  -1.93%  "

# Row 37 - dogwifhat
$ws.Range("E37").Value = "  +3.35%  "

# Row 38 - Hedera
$ws.Range("E38").Value = "  -0.79%  "

# Row 39: was Mantle -> now Kaspa (rows 39/40 swap order, with refreshed data)
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D39" "0.137"
$ws.Range("E39").Value = "  -1.78%  "

# Row 40: was Kaspa -> now Mantle
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D40" "0.999"
$ws.Range("E40").Value = "  -1.85%  "

# Row 41 - Filecoin
Set-TextValue "D41" "5.84"
$ws.Range("E41").Value = "  -1.67%  "

# Row 42 - FirstDigitalUSD
$ws.Range("E42").Value = "  +0.14%  "

# Row 43 - TheGraph
$ws.Range("E43").Value = "  -1.13%  "

# Row 45 - Cosmos
Set-TextValue "D45" "8.73"
$ws.Range("E45").Value = "  -0.41%  "

# Row 46 - Stacks
$ws.Range("E46").Value = "  -2.42%  "

# Row 47 - OKB
Set-TextValue "D47" "45.92"
$ws.Range("E47").Value = "  -2.51%  "

# Row 48 - Bittensor
Set-TextValue "D48" "397.56"
$ws.Range("E48").Value = "  -4.74%  "

# Row 49 - FLOKI
$ws.Range("E49").Value = "  -8.97%  "

# Row 50 - VeChain
$ws.Range("E50").Value = "  -1.59%  "

# Row 51: Monero -> Arweave
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D51" "38.96"
$ws.Range("E51").Value = "  +2.04%  "
